# Generate Report for Handback
# The 9a6476c5-... file has finished translation and is now handed back.
# Update its Status to "Handed back: in sync with en-US" and record the
# new handback timestamps on the per-language sheets, everywhere this
# row is reported: the Overview summary sheet and each language sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-20 06:37:30"

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-20 06:37:35"
